# The "statut" column used a black/red/orange/green emoji-square legend;
# this recolors the "black" status to "blue" (emoji + label), leaving the
# red/orange/green entries untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange

# LookAt:=1 -> xlWhole, so only cells whose *entire* content equals the
# search text are changed (avoids accidentally touching substrings in
# unrelated longer strings).
[void]$ur.Replace("⬛", "📘", 1)
[void]$ur.Replace("🟥", "📕", 1)
[void]$ur.Replace("🟧", "📙", 1)
[void]$ur.Replace("🟩", "📗", 1)
[void]$ur.Replace("noir", "bleu", 1)
